# Apply the "Add files via upload" edit:
#   - Product column (I) renamed from "Reputation Booster" to "Review Generator"
#     for every data row.
#   - A handful of Email/Phone values were edited down to truncated/partial
#     text (as if the user started retyping them), and Excel created brand
#     new mailto hyperlinks for the cells that didn't already have one
#     (A2, A12, A14) while leaving the pre-existing hyperlink on A5 intact.
#   - Column I was narrowed to fit the shorter "Review Generator" text.
#   - The active selection ended up on D2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Product column: "Reputation Booster" -> "Review Generator" ----------
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 9).Value = "Review Generator"
}

# --- Edited Email / Phone cells -------------------------------------------
# (order matches the original authoring session so newly-introduced shared
# strings land at the same indices as the target file)

# Row 5 (Ella Ballard) -- already has a hyperlink, keep it, just edit text
$ws.Range("A5").Value = "EllaBallard@gmai"

# Row 12 (Emmett Baker)
$ws.Range("A12").Value = "EmmettBaker@gm"

# Row 14 (Sydney Macy)
$ws.Range("B14").Value = "(865) 519-46"
$ws.Range("A14").Value = "SydneyMacy@gma"

# Row 5 phone
$ws.Range("B5").Value = "(652) 771"

# Row 2 (Scarlett Titterington)
$ws.Range("A2").Value = "ScarlettTitterington@"
$ws.Range("B2").Value = "(651) 9"

# --- New hyperlinks for the cells that didn't already have one -----------
# Order matters for relationship-id allocation: A12, then A14, then A2.
# (Hyperlinks.Add re-styles the cell with a fresh "Hyperlink" xf; the column
# was already formatted that way via style index carried on A3, so restore
# it to keep the cell's style identical to its neighbours / to before.)
$ws.Hyperlinks.Add($ws.Range("A12"), "mailto:EmmettBaker@gm")
$ws.Range("A12").Style = $ws.Range("A3").Style

$ws.Hyperlinks.Add($ws.Range("A14"), "mailto:SydneyMacy@gma")
$ws.Range("A14").Style = $ws.Range("A3").Style

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:ScarlettTitterington@")
$ws.Range("A2").Style = $ws.Range("A3").Style

# --- Column I width: narrower now that "Review Generator" is shorter -----
$ws.Columns.Item(9).ColumnWidth = 16.14

# --- Final selection on D2 --------------------------------------------------
$ws.Range("D2").Select()
